# Auto-generated: fix page-num bug by appending rows 22-48 to pageList sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Range("A22").Value = "tag-ill"
$ws.Range("B22").Value = "這是第五篇文章21"
$ws.Range("C22").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D22").Value = "這是第六篇文章的描述"
$ws.Range("E22").Value = "img/index/test3.jpg"
$ws.Range("F22").Value = "2014/3/24"

# Row 23
$ws.Range("A23").Value = "tag-css"
$ws.Range("B23").Value = "這是第六篇文章22"
$ws.Range("C23").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D23").Value = "這是第七篇文章的描述"
$ws.Range("E23").Value = "img/index/test1.jpg"
$ws.Range("F23").Value = "2014/2/24"

# Row 24
$ws.Range("A24").Value = "tag-ui"
$ws.Range("B24").Value = "這是第五篇文章23"
$ws.Range("C24").Value = "http://tw.yahoo.com"
$ws.Range("D24").Value = "這是第一篇文章的描述"
$ws.Range("E24").Value = "img/index/test1.jpg"
$ws.Range("F24").Value = "2014/3/25"

# Row 25
$ws.Range("A25").Value = "tag-ill"
$ws.Range("B25").Value = "這是第六篇文章24"
$ws.Range("C25").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D25").Value = "這是第二篇文章的描述"
$ws.Range("E25").Value = "img/index/test2.jpg"
$ws.Range("F25").Value = "2014/2/25"

# Row 26
$ws.Range("A26").Value = "tag-design"
$ws.Range("B26").Value = "這是第五篇文章25"
$ws.Range("C26").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D26").Value = "這是第三篇文章的描述"
$ws.Range("E26").Value = "img/index/test3.jpg"
$ws.Range("F26").Value = "2014/3/26"

# Row 27
$ws.Range("A27").Value = "tag-photo"
$ws.Range("B27").Value = "這是第六篇文章26"
$ws.Range("C27").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D27").Value = "這是第四篇文章的描描述描述描述描述描述描述描述描述描述描述描述描述描述述描述描述描述描述描述描述描述描述"
$ws.Range("E27").Value = "img/index/test1.jpg"
$ws.Range("F27").Value = "2014/2/26"

# Row 28
$ws.Range("A28").Value = "tag-web"
$ws.Range("B28").Value = "這是第五篇文章27"
$ws.Range("C28").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D28").Value = "這是第五篇文章的描述"
$ws.Range("E28").Value = "img/index/test2.jpg"
$ws.Range("F28").Value = "2014/3/27"

# Row 29
$ws.Range("A29").Value = "tag-others"
$ws.Range("B29").Value = "這是第六篇文章28"
$ws.Range("C29").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D29").Value = "這是第六篇文章的描述"
$ws.Range("E29").Value = "img/index/test3.jpg"
$ws.Range("F29").Value = "2014/2/27"

# Row 30
$ws.Range("A30").Value = "tag-ill"
$ws.Range("B30").Value = "這是第篇文章28"
$ws.Range("C30").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D30").Value = "這是第六篇文章的描述"
$ws.Range("E30").Value = "img/index/test3.jpg"
$ws.Range("F30").Value = "2014/2/27"

# Row 31
$ws.Range("A31").Value = "tag-others"
$ws.Range("B31").Value = "這是第篇文章29"
$ws.Range("C31").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D31").Value = "這是第六篇文章的描述"
$ws.Range("E31").Value = "img/index/test4.jpg"
$ws.Range("F31").Value = "2014/2/28"

# Row 32
$ws.Range("A32").Value = "tag-design"
$ws.Range("B32").Value = "這是第篇文章30"
$ws.Range("C32").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D32").Value = "這是第六篇文章的描述"
$ws.Range("E32").Value = "img/index/test5.jpg"
$ws.Range("F32").Value = "2014/2/29"

# Row 33
$ws.Range("A33").Value = "tag-photo"
$ws.Range("B33").Value = "這是第篇文章31"
$ws.Range("C33").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D33").Value = "這是第六篇文章的描述"
$ws.Range("E33").Value = "img/index/test6.jpg"
$ws.Range("F33").Value = "2014/2/30"

# Row 34
$ws.Range("A34").Value = "tag-web"
$ws.Range("B34").Value = "這是第篇文章32"
$ws.Range("C34").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D34").Value = "這是第六篇文章的描述"
$ws.Range("E34").Value = "img/index/test1.jpg"
$ws.Range("F34").Value = "2014/2/31"

# Row 35
$ws.Range("A35").Value = "tag-ill"
$ws.Range("B35").Value = "這是第篇文章33"
$ws.Range("C35").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D35").Value = "這是第六篇文章的描述"
$ws.Range("E35").Value = "img/index/test2.jpg"
$ws.Range("F35").Value = "2014/2/32"

# Row 36
$ws.Range("A36").Value = "tag-ill"
$ws.Range("B36").Value = "這是第篇文章34"
$ws.Range("C36").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D36").Value = "這是第六篇文章的描述"
$ws.Range("E36").Value = "img/index/test3.jpg"
$ws.Range("F36").Value = "2014/2/33"

# Row 37
$ws.Range("A37").Value = "tag-css"
$ws.Range("B37").Value = "這是第篇文章35"
$ws.Range("C37").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D37").Value = "這是第六篇文章的描述"
$ws.Range("E37").Value = "img/index/test3.jpg"
$ws.Range("F37").Value = "2014/2/34"

# Row 38
$ws.Range("A38").Value = "tag-ui"
$ws.Range("B38").Value = "這是第篇文章36"
$ws.Range("C38").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D38").Value = "這是第六篇文章的描述"
$ws.Range("E38").Value = "img/index/test4.jpg"
$ws.Range("F38").Value = "2014/2/35"

# Row 39
$ws.Range("A39").Value = "tag-ill"
$ws.Range("B39").Value = "這是第篇文章37"
$ws.Range("C39").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D39").Value = "這是第六篇文章的描述"
$ws.Range("E39").Value = "img/index/test5.jpg"
$ws.Range("F39").Value = "2014/2/36"

# Row 40
$ws.Range("A40").Value = "tag-design"
$ws.Range("B40").Value = "這是第篇文章38"
$ws.Range("C40").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D40").Value = "這是第六篇文章的描述"
$ws.Range("E40").Value = "img/index/test6.jpg"
$ws.Range("F40").Value = "2014/2/37"

# Row 41
$ws.Range("A41").Value = "tag-photo"
$ws.Range("B41").Value = "這是第篇文章39"
$ws.Range("C41").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D41").Value = "這是第六篇文章的描述"
$ws.Range("E41").Value = "img/index/test1.jpg"
$ws.Range("F41").Value = "2014/2/38"

# Row 42
$ws.Range("A42").Value = "tag-web"
$ws.Range("B42").Value = "這是第篇文章40"
$ws.Range("C42").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D42").Value = "這是第六篇文章的描述"
$ws.Range("E42").Value = "img/index/test2.jpg"
$ws.Range("F42").Value = "2014/2/39"

# Row 43
$ws.Range("A43").Value = "tag-others"
$ws.Range("B43").Value = "這是第篇文章41"
$ws.Range("C43").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D43").Value = "這是第六篇文章的描述"
$ws.Range("E43").Value = "img/index/test3.jpg"
$ws.Range("F43").Value = "2014/2/40"

# Row 44
$ws.Range("A44").Value = "tag-ill"
$ws.Range("B44").Value = "這是第篇文章42"
$ws.Range("C44").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D44").Value = "這是第六篇文章的描述"
$ws.Range("E44").Value = "img/index/test3.jpg"
$ws.Range("F44").Value = "2014/2/41"

# Row 45
$ws.Range("A45").Value = "tag-photo"
$ws.Range("B45").Value = "這是第篇文章43"
$ws.Range("C45").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D45").Value = "這是第六篇文章的描述"
$ws.Range("E45").Value = "img/index/test4.jpg"
$ws.Range("F45").Value = "2014/2/42"

# Row 46
$ws.Range("A46").Value = "tag-web"
$ws.Range("B46").Value = "這是第篇文章44"
$ws.Range("C46").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D46").Value = "這是第六篇文章的描述"
$ws.Range("E46").Value = "img/index/test5.jpg"
$ws.Range("F46").Value = "2014/2/43"

# Row 47
$ws.Range("A47").Value = "tag-others"
$ws.Range("B47").Value = "這是第篇文章45"
$ws.Range("C47").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D47").Value = "這是第六篇文章的描述"
$ws.Range("E47").Value = "img/index/test6.jpg"
$ws.Range("F47").Value = "2014/2/44"

# Row 48
$ws.Range("A48").Value = "tag-ill"
$ws.Range("B48").Value = "這是第篇文章46"
$ws.Range("C48").Value = "http://oxxo-studio.blogspot.com"
$ws.Range("D48").Value = "這是第六篇文章的描述"
$ws.Range("E48").Value = "img/index/test6.jpg"
$ws.Range("F48").Value = "2014/2/45"

# Update the view: zoom 85%, scroll so row 16 is at top, select I48 (matches the post-edit sheetView)
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("I48").Select()

